# Weekly update: insert a new week's price record as the new top data row
# (row 174) for the "Pepino ensalada" sheet, shifting all subsequent rows
# down by one (old row 174 -> 175, ..., old row 207 -> 208).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 174..207 down to 175..208, leaving row 174 free.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with this week's record.
$ws.Range("A174").Value = 7
$ws.Range("B174").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C174").Value = "Ñuble"
$ws.Range("D174").Value = 44617
$ws.Range("E174").Value = 16
$ws.Range("F174").Value = 100112043
$ws.Range("G174").Value = "Pepino ensalada"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 120
$ws.Range("K174").Value = 13000
$ws.Range("L174").Value = 13500
$ws.Range("M174").Value = 13250
$ws.Range("N174").Value = "$/caja 80 unidades"
$ws.Range("O174").Value = "Región del Maule"
$ws.Range("P174").Value = 166
$ws.Range("Q174").Value = 80
$ws.Range("R174").Value = "Hortaliza"
